$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (rows 2-9) holds the "Förändrad" date, stored as a serial date number.
# Bump each of these dates forward by one day (46074 -> 46075).
for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value = 46075
    }
}
